# Reorder rows 2-37 on the active sheet per the target permutation.
# Read every source row (A:T) into memory first via Value2 (plain property,
# unlike the parameterized Value getter) so writes never clobber an
# as-yet-unread source row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 20
$firstRow = 2
$lastRow = 37

$original = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
  $rowVals = @()
  for ($c = 1; $c -le $lastCol; $c++) {
    $rowVals += ,$ws.Cells.Item($r, $c).Value2
  }
  $original[$r] = $rowVals
}

# target row -> source row (source row data as it existed before this edit)
$targetToSource = @{
  2 = 29
  3 = 15
  4 = 9
  5 = 6
  6 = 7
  7 = 31
  8 = 34
  9 = 20
  10 = 26
  11 = 10
  12 = 27
  13 = 5
  14 = 28
  15 = 18
  16 = 25
  17 = 35
  18 = 30
  19 = 36
  20 = 23
  21 = 17
  22 = 19
  23 = 3
  24 = 21
  25 = 11
  26 = 22
  27 = 32
  28 = 33
  29 = 37
  30 = 13
  31 = 14
  32 = 16
  33 = 12
  34 = 2
  35 = 8
  36 = 24
  37 = 4
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
  $src = $targetToSource[$r]
  $srcVals = $original[$src]
  for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item($r, $c).Value = $srcVals[$c-1]
  }
}

Write-Output "Reordered rows $firstRow..$lastRow"